$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 81; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "BaseLoader(resize=(128, 128))") {
        $cell.Value2 = "SimpleLoader(resize=(128, 128))"
    }
}
